# "update to the dark mode" -- Pendientes.xlsx status-tracker update:
#   - Backend "Agregar productos con varias imagenes" (B4) -> OK (green)
#   - Backend "Definir acceso con roles..." (B5) -> OK (green)
#   - Frontend row for "Definir acceso con roles..." task gets filled in
#     (D3 description + E3 status = OK, green) in the Frontend block
#   - Column D widened, a thin spacer column G introduced
#   - Selection moved to D13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$green = 5287936  # RGB(00,176,80) == fgColor FF00B050, BGR-packed OLE color

# --- Frontend block, row 3: add the missing task description + status ---
$d3 = $ws.Range("D3")
$d3.Value = "Definir acceso con roles (Admin, cliente, vendedor)"

# Clone B3's "OK"/green format onto E3 in a single style transition (avoids
# leaving an orphaned intermediate cellXf behind), then set its value.
$ws.Range("B3").Copy()
$ws.Range("E3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$e3 = $ws.Range("E3")
$e3.Value = "OK"

# --- Backend block: both pending tasks flip to OK / green ---
$b4 = $ws.Range("B4")
$b4.Value = "OK"
$b4.Interior.Color = $green

$b5 = $ws.Range("B5")
$b5.Value = "OK"
$b5.Interior.Color = $green

# --- column widths: widen Description (D), add thin spacer column (G) ---
$ws.Columns.Item(4).ColumnWidth = 60.666666666666664
$ws.Columns.Item(7).ColumnWidth = 2.6666666666666665

# --- move the active selection ---
$ws.Range("D13").Select()
